$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.89

$ws.Range("F3").Value = 3.5
$ws.Range("I3").Value = 2.18
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 2.32
$ws.Range("P3").Value = 2.02
$ws.Range("S3").Value = 2.44
$ws.Range("T3").Value = 1.52
$ws.Range("U3").Value = 1.04
$ws.Range("V3").Value = 1.84
$ws.Range("X3").Value = 990
$ws.Range("Y3").Value = 990
$ws.Range("AB3").Value = 990
$ws.Range("AC3").Value = 990
$ws.Range("AD3").Value = 990
$ws.Range("AG3").Value = 990
$ws.Range("AH3").Value = 990

$ws.Range("J4").Value = 1.03
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 1.31
$ws.Range("O4").Value = 1.01
$ws.Range("P4").Value = 1.3
$ws.Range("R4").Value = 1.13
$ws.Range("S4").Value = 1.05
$ws.Range("T4").Value = 1.04
$ws.Range("U4").Value = 1.04
$ws.Range("V4").Value = 1.01
$ws.Range("W4").Value = 1.01
$ws.Range("X4").Value = 990
$ws.Range("Y4").Value = 990
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 990
$ws.Range("AC4").Value = 990
$ws.Range("AD4").Value = 990
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 990
$ws.Range("AH4").Value = 990
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

$ws.Range("G5").Value = 1.94
$ws.Range("I5").Value = 5.2
$ws.Range("J5").Value = 3.8
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 2.2
$ws.Range("O5").Value = 1.26
$ws.Range("P5").Value = 2.04
$ws.Range("R5").Value = 1.33
$ws.Range("S5").Value = 2.64
$ws.Range("T5").Value = 1.58
$ws.Range("U5").Value = 1.04
$ws.Range("V5").Value = 1.23
$ws.Range("W5").Value = 2.06
$ws.Range("X5").Value = 990
$ws.Range("Y5").Value = 990
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 990
$ws.Range("AC5").Value = 990
$ws.Range("AD5").Value = 990
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 990
$ws.Range("AH5").Value = 990
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000

$ws.Range("F6").Value = 4.7
$ws.Range("G6").Value = 5.7
$ws.Range("H6").Value = 1.58
$ws.Range("I6").Value = 1.72
$ws.Range("J6").Value = 4.6
$ws.Range("K6").Value = 5.6
$ws.Range("P6").Value = 2.68
$ws.Range("Q6").Value = 1.47

$ws.Range("P7").Value = 1.92
$ws.Range("Q7").Value = 2.02
$ws.Range("R7").Value = 1.35

$ws.Range("G8").Value = 7.2
$ws.Range("H8").Value = 1.58
$ws.Range("I8").Value = 1.59
$ws.Range("O8").Value = 1.33
$ws.Range("P8").Value = 1.96
$ws.Range("Q8").Value = 1.99
$ws.Range("R8").Value = 1.36
$ws.Range("S8").Value = 3.55
$ws.Range("T8").Value = 2.06
$ws.Range("X8").Value = 15.5
$ws.Range("Z8").Value = 8.800000000000001
$ws.Range("AB8").Value = 20
$ws.Range("AC8").Value = 9.6
$ws.Range("AF8").Value = 60
$ws.Range("AG8").Value = 26
$ws.Range("AJ8").Value = 260
$ws.Range("AL8").Value = 120
$ws.Range("AN8").Value = 180
$ws.Range("AO8").Value = 9.4

$ws.Range("H9").Value = 4.4
$ws.Range("N9").Value = 2.32
$ws.Range("O9").Value = 1.75
$ws.Range("P9").Value = 1.4
